$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Label" header column (H): binary Control(0) / MDD(1) flag per patient row ---
$ws.Range("H1").Value = "Label"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats - match the other header cells' style

# --- Refit results: updated Prediction (D) / Error (E) values for the 100-iteration block ---
$ws.Range("D2").Value = 0.6245256226085666
$ws.Range("E2").Value = 0.6245256226085666

$ws.Range("D5").Value = 0.07279741489621196
$ws.Range("E5").Value = 0.07279741489621196

$ws.Range("D6").Value = 0.604711681971891
$ws.Range("E6").Value = 0.604711681971891

$ws.Range("D7").Value = 0.3152016383198216
$ws.Range("E7").Value = 0.6847983616801785

$ws.Range("D8").Value = 0.6151953057175282
$ws.Range("E8").Value = 0.3848046942824718

$ws.Range("D9").Value = 0.4424813718845387
$ws.Range("E9").Value = 0.5575186281154614

$ws.Range("D10").Value = 0.4940514275738653
$ws.Range("E10").Value = 0.5059485724261348

# --- New Label column values: 0 = Control, 1 = MDD (both the 100- and 200-iteration blocks) ---
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("H7").Value = 1
$ws.Range("H8").Value = 1
$ws.Range("H9").Value = 1
$ws.Range("H10").Value = 1
$ws.Range("H11").Value = 1

$ws.Range("H12").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("H17").Value = 1
$ws.Range("H18").Value = 1
$ws.Range("H19").Value = 1
$ws.Range("H20").Value = 1
$ws.Range("H21").Value = 1
